$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (shifts old row17->18, old row18->19)
$ws.Rows("17").Insert()

# Populate new row 17 as a duplicate of row 16 (same formatting/content),
# then bump the sequence number in column A to 11.
$ws.Range("A16:Q16").Copy($ws.Range("A17:Q17"))
$ws.Rows("17").RowHeight = 25.5
$ws.Range("A17").Value = 11

# Update the totals row (now row 18) to reflect the new item's price addition.
$ws.Range("N18").Value = 398.48500000000001

Write-Host "done"
